# CCC19_Derived_Variables_Spreadsheet.xlsx - add "Metastatic cancer status"
# derived variable row (Ca19 / metastatic / Cancer) to Table1, inserted as
# the new row 50 (pushing the existing rows 50-177 down to 51-178).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new worksheet row above row 50 (shifts everything below it down
# by one row, including the remaining table rows).
$ws.Rows.Item(50).Insert()

# Populate the new row with the new derived-variable metadata.
$ws.Range("A50").Value = "Ca19"
$ws.Range("B50").Value = "metastatic"
$ws.Range("C50").Value = "Cancer"
$ws.Range("D50").Value = "Metastatic cancer status"

# Grow Table1 so the new row (and the now one-row-taller data range) is
# included again, keeping the autofilter/table ref in sync.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E178"))

# Match the author's final selection/scroll position in the saved file.
$ws.Range("B50").Select()
